# Auto-generated edit script: updates cryptos list cell values per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to remain a literal text string (matches the source
    # data, which stores every value as inlineStr/shared-string text even
    # when it looks numeric, e.g. "1.0000" or "26.258.77").
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "26.258.77"
Set-TextCell "E2" "  +2.90%  "
Set-TextCell "D3" "1.719.80"
Set-TextCell "E3" "  +3.38%  "
Set-TextCell "E4" "  +0.05%  "
Set-TextCell "D5" "239.69"
Set-TextCell "E5" "  +1.10%  "
Set-TextCell "D6" "1.0000"
Set-TextCell "E6" "  -0.01%  "
Set-TextCell "D7" "0.4710"
Set-TextCell "E7" "  -1.84%  "
Set-TextCell "E8" "  -0.04%  "
Set-TextCell "E9" "  +0.69%  "
Set-TextCell "D10" "1.715.46"
Set-TextCell "E10" "  +3.22%  "
Set-TextCell "D11" "0.07068"
Set-TextCell "E11" "  -0.54%  "
Set-TextCell "D12" "15.34"
Set-TextCell "D13" "0.5973"
Set-TextCell "E13" "  +1.67%  "
Set-TextCell "E14" "  +1.22%  "
Set-TextCell "E15" "  +1.91%  "
Set-TextCell "D16" "1.0000"
Set-TextCell "E16" "  -0.03%  "
Set-TextCell "D17" "1.000"
Set-TextCell "E17" "  +0.04%  "
Set-TextCell "D18" "26.270.61"
Set-TextCell "E18" "  +3.01%  "
Set-TextCell "D19" "0.000006802"
Set-TextCell "E19" "  +0.96%  "
Set-TextCell "D20" "11.55"
Set-TextCell "E20" "  +0.79%  "
Set-TextCell "D21" "1.935.22"
Set-TextCell "E21" "  +3.33%  "
Set-TextCell "D22" "4.543"
Set-TextCell "E22" "  +2.56%  "
Set-TextCell "D23" "8.724"
Set-TextCell "E23" "  +0.61%  "
Set-TextCell "D24" "5.279"
Set-TextCell "E24" "  +0.14%  "
Set-TextCell "D25" "134.78"
Set-TextCell "E25" "  +0.82%  "
Set-TextCell "D26" "15.23"
Set-TextCell "E26" "  +1.22%  "
Set-TextCell "D27" "1.401"
Set-TextCell "E27" "  +0.96%  "
Set-TextCell "D28" "1.760"
Set-TextCell "E28" "  +2.68%  "
Set-TextCell "D29" "107.43"
Set-TextCell "E29" "  +1.98%  "
Set-TextCell "D30" "3.977"
Set-TextCell "E30" "  +0.69%  "
Set-TextCell "D31" "3.679"
Set-TextCell "E31" "  +0.33%  "
Set-TextCell "D32" "0.07762"
Set-TextCell "E32" "  +1.36%  "
Set-TextCell "D33" "0.04452"
Set-TextCell "E33" "  +5.77%  "
Set-TextCell "D34" "2.614"
Set-TextCell "E34" "  +0.16%  "
Set-TextCell "D35" "0.9757"
Set-TextCell "E35" "  +2.64%  "
Set-TextCell "D36" "0.6178"
Set-TextCell "E36" "  +1.29%  "
Set-TextCell "D37" "0.9323"
Set-TextCell "E37" "  +7.47%  "
Set-TextCell "D38" "111.89"
Set-TextCell "E38" "  +15.60%  "
Set-TextCell "D39" "2.416"
Set-TextCell "E39" "  -7.09%  "
Set-TextCell "D40" "1.925"
Set-TextCell "E40" "  +3.80%  "
Set-TextCell "D41" "0.9998"
Set-TextCell "D42" "0.01480"
Set-TextCell "E42" "  +0.88%  "
Set-TextCell "D43" "5.419"
Set-TextCell "E43" "  +12.91%  "
Set-TextCell "D44" "0.3820"
Set-TextCell "E44" "  +1.57%  "
Set-TextCell "D45" "0.1180"
Set-TextCell "E45" "  +4.80%  "
Set-TextCell "E46" "  +0.92%  "
Set-TextCell "D47" "0.05265"
Set-TextCell "E47" "  +0.25%  "
Set-TextCell "B48" "EnergySwap"
Set-TextCell "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D48" "7.783"
Set-TextCell "E48" "  +6.78%  "
Set-TextCell "B49" "Elrond"
Set-TextCell "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell "D49" "30.17"
Set-TextCell "E49" "  +1.69%  "
Set-TextCell "D50" "0.3380"
Set-TextCell "E50" "  +1.77%  "
Set-TextCell "B51" "Aave"
Set-TextCell "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D51" "50.69"
Set-TextCell "E51" "  +1.79%  "
